$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell D1, matching the header style used by A1:C1 (bold, bordered,
# centered/top-aligned) by copying the existing header cell's formatting.
$ws.Range("D1").Value = "Ano"
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Fill D2:D10 with the reference year value
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = 2023
}
